# Add two new BOM line items (Voltage Regulator, Slide Switch) under the
# existing "ATMega328P-PU" row, each linking out to the newly added
# gerber / eagle design files, and update the saved cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM rows
$ws.Range("C6").Value = "Voltage Regulator"
$ws.Range("C7").Value = "Slide Switch"

# Hyperlink the new items out to the newly committed gerber / eagle files
$c6 = $ws.Range("C6")
$c6.Hyperlinks.Add($c6, "https://github.com/nicholasnguyen/micromouse/blob/master/gerber.zip")

$c7 = $ws.Range("C7")
$c7.Hyperlinks.Add($c7, "https://github.com/nicholasnguyen/micromouse/blob/master/eagle.zip")

# Adding the hyperlink re-applies formatting; reset the cells back to the
# workbook's shared "Hyperlink" cell style (same style already used by the
# other BOM hyperlinks such as C9:C12).
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("C7").Style = "Hyperlink"

# Move / save the active selection as it was left after the edit
[void]$ws.Range("C17").Select()
